$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-25 07:29:30"
$wsZhCn.Range("H4").Value = "2016-03-25 07:29:59"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-25 07:29:34"
$wsDeDe.Range("H4").Value = "2016-03-25 07:30:10"
